$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 130960607
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("M5").Value = 'äldre spår'
$ws.Range("Q5").Value = 446240
$ws.Range("R5").Value = 6759818
$ws.Range("Z5").Value = '10:26'
$ws.Range("AB5").Value = '10:26'
$ws.Range("AC5").Value = ''
$ws.Range("A6").Value = 130963816
$ws.Range("B6").Value = 79244
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = 'Garnlav'
$ws.Range("G6").Value = 'Alectoria sarmentosa'
$ws.Range("H6").Value = '(Ach.) Ach.'
$ws.Range("M6").Value = ''
$ws.Range("Q6").Value = 445932
$ws.Range("R6").Value = 6760079
$ws.Range("Z6").Value = '14:08'
$ws.Range("AB6").Value = '14:08'
$ws.Range("AC6").Value = 'Rikligt i närområdet'
$ws.Range("A7").Value = 130961458
$ws.Range("Q7").Value = 446059
$ws.Range("R7").Value = 6760088
$ws.Range("A8").Value = 130961962
$ws.Range("Q8").Value = 446084
$ws.Range("R8").Value = 6759981
$ws.Range("A15").Value = 130963950
$ws.Range("Q15").Value = 445926
$ws.Range("R15").Value = 6760113
$ws.Range("Z15").Value = '14:08'
$ws.Range("AB15").Value = '14:08'
$ws.Range("A16").Value = 130961105
$ws.Range("Q16").Value = 446124
$ws.Range("R16").Value = 6759989
$ws.Range("Z16").Value = '10:26'
$ws.Range("AB16").Value = '10:26'
$ws.Range("A19").Value = 130961956
$ws.Range("B19").Value = 79863
$ws.Range("E19").Value = 6453
$ws.Range("F19").Value = 'Vedskivlav'
$ws.Range("G19").Value = 'Hertelidea botryosa'
$ws.Range("H19").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("M19").Value = ''
$ws.Range("Q19").Value = 446084
$ws.Range("R19").Value = 6759981
$ws.Range("AC19").Value = 'Miljöbilder'
$ws.Range("A20").Value = 130960378
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = 'Tretåig hackspett'
$ws.Range("G20").Value = 'Picoides tridactylus'
$ws.Range("H20").Value = '(Linnaeus, 1758)'
$ws.Range("M20").Value = 'äldre spår'
$ws.Range("Q20").Value = 446272
$ws.Range("R20").Value = 6759739
$ws.Range("AC20").Value = ''
$ws.Range("A25").Value = 130962090
$ws.Range("B25").Value = 79244
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = 'Garnlav'
$ws.Range("G25").Value = 'Alectoria sarmentosa'
$ws.Range("H25").Value = '(Ach.) Ach.'
$ws.Range("M25").Value = ''
$ws.Range("Q25").Value = 446080
$ws.Range("R25").Value = 6759960
$ws.Range("A26").Value = 130961746
$ws.Range("B26").Value = 57881
$ws.Range("E26").Value = 100049
$ws.Range("F26").Value = 'Spillkråka'
$ws.Range("G26").Value = 'Dryocopus martius'
$ws.Range("H26").Value = '(Linnaeus, 1758)'
$ws.Range("M26").Value = 'färska spår'
$ws.Range("Q26").Value = 446098
$ws.Range("R26").Value = 6760061
$ws.Range("A27").Value = 130961461
$ws.Range("B27").Value = 79244
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = 'Garnlav'
$ws.Range("G27").Value = 'Alectoria sarmentosa'
$ws.Range("H27").Value = '(Ach.) Ach.'
$ws.Range("M27").Value = ''
$ws.Range("Q27").Value = 446088
$ws.Range("R27").Value = 6760088
$ws.Range("Z27").Value = '10:26'
$ws.Range("AB27").Value = '10:26'
$ws.Range("A28").Value = 130961750
$ws.Range("Q28").Value = 446098
$ws.Range("R28").Value = 6760061
$ws.Range("AC28").Value = 'Rikligt i en radie av ca 50 meter'
$ws.Range("A29").Value = 130963807
$ws.Range("B29").Value = 57881
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = 'Spillkråka'
$ws.Range("G29").Value = 'Dryocopus martius'
$ws.Range("H29").Value = '(Linnaeus, 1758)'
$ws.Range("M29").Value = 'färska spår'
$ws.Range("Q29").Value = 445932
$ws.Range("R29").Value = 6760079
$ws.Range("Z29").Value = '14:08'
$ws.Range("AB29").Value = '14:08'
$ws.Range("AC29").Value = ''
